{"js": "// Update the division-problem answers in the table to the new set of\n// three-digit \u00f7 one-digit practice problems.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"237\u00f78=29, 5\", \"656\u00f73=218, 2\"],\n  [\"612\u00f75=122, 2\", \"208\u00f75=41, 3\"],\n  [\"972\u00f77=138, 6\", \"658\u00f73=219, 1\"],\n  [\"916\u00f74=229, 0\", \"191\u00f77=27, 2\"],\n  [\"555\u00f75=111, 0\", \"757\u00f72=378, 1\"],\n  [\"585\u00f74=146, 1\", \"626\u00f74=156, 2\"],\n  [\"984\u00f77=140, 4\", \"230\u00f78=28, 6\"],\n  [\"120\u00f72=60, 0\", \"242\u00f73=80, 2\"],\n  [\"455\u00f77=65, 0\", \"738\u00f73=246, 0\"],\n  [\"423\u00f72=211, 1\", \"523\u00f76=87, 1\"],\n  [\"939\u00f74=234, 3\", \"483\u00f75=96, 3\"],\n  [\"585\u00f76=97, 3\", \"152\u00f74=38, 0\"],\n  [\"734\u00f73=244, 2\", \"943\u00f76=157, 1\"],\n  [\"202\u00f77=28, 6\", \"892\u00f72=446, 0\"],\n  [\"490\u00f73=163, 1\", \"958\u00f75=191, 3\"],\n  [\"230\u00f75=46, 0\", \"871\u00f75=174, 1\"],\n  [\"313\u00f78=39, 1\", \"216\u00f73=72, 0\"],\n  [\"733\u00f74=183, 1\", \"659\u00f74=164, 3\"],\n  [\"929\u00f74=232, 1\", \"950\u00f74=237, 2\"],\n  [\"177\u00f75=35, 2\", \"485\u00f73=161, 2\"],\n  [\"259\u00f78=32, 3\", \"569\u00f72=284, 1\"],\n  [\"512\u00f76=85, 2\", \"786\u00f72=393, 0\"],\n  [\"772\u00f76=128, 4\", \"134\u00f75=26, 4\"],\n  [\"355\u00f72=177, 1\", \"151\u00f75=30, 1\"],\n  [\"356\u00f79=39, 5\", \"587\u00f76=97, 5\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem answers in the table to the new set of\n# three-digit \u00f7 one-digit practice problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"237\u00f78=29, 5\";  New = \"656\u00f73=218, 2\" },\n    @{ Old = \"612\u00f75=122, 2\"; New = \"208\u00f75=41, 3\" },\n    @{ Old = \"972\u00f77=138, 6\"; New = \"658\u00f73=219, 1\" },\n    @{ Old = \"916\u00f74=229, 0\"; New = \"191\u00f77=27, 2\" },\n    @{ Old = \"555\u00f75=111, 0\"; New = \"757\u00f72=378, 1\" },\n    @{ Old = \"585\u00f74=146, 1\"; New = \"626\u00f74=156, 2\" },\n    @{ Old = \"984\u00f77=140, 4\"; New = \"230\u00f78=28, 6\" },\n    @{ Old = \"120\u00f72=60, 0\";  New = \"242\u00f73=80, 2\" },\n    @{ Old = \"455\u00f77=65, 0\";  New = \"738\u00f73=246, 0\" },\n    @{ Old = \"423\u00f72=211, 1\"; New = \"523\u00f76=87, 1\" },\n    @{ Old = \"939\u00f74=234, 3\"; New = \"483\u00f75=96, 3\" },\n    @{ Old = \"585\u00f76=97, 3\";  New = \"152\u00f74=38, 0\" },\n    @{ Old = \"734\u00f73=244, 2\"; New = \"943\u00f76=157, 1\" },\n    @{ Old = \"202\u00f77=28, 6\";  New = \"892\u00f72=446, 0\" },\n    @{ Old = \"490\u00f73=163, 1\"; New = \"958\u00f75=191, 3\" },\n    @{ Old = \"230\u00f75=46, 0\";  New = \"871\u00f75=174, 1\" },\n    @{ Old = \"313\u00f78=39, 1\";  New = \"216\u00f73=72, 0\" },\n    @{ Old = \"733\u00f74=183, 1\"; New = \"659\u00f74=164, 3\" },\n    @{ Old = \"929\u00f74=232, 1\"; New = \"950\u00f74=237, 2\" },\n    @{ Old = \"177\u00f75=35, 2\";  New = \"485\u00f73=161, 2\" },\n    @{ Old = \"259\u00f78=32, 3\";  New = \"569\u00f72=284, 1\" },\n    @{ Old = \"512\u00f76=85, 2\";  New = \"786\u00f72=393, 0\" },\n    @{ Old = \"772\u00f76=128, 4\"; New = \"134\u00f75=26, 4\" },\n    @{ Old = \"355\u00f72=177, 1\"; New = \"151\u00f75=30, 1\" },\n    @{ Old = \"356\u00f79=39, 5\";  New = \"587\u00f76=97, 5\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute([ref]$pair.Old, [ref]$false, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$pair.New, [ref]2) | Out-Null\n}\n\n$d.Save()\n"}
